# Deploy the implementation guide.
#
# This reproduces the OOXML changes:
#   - sharedStrings: "active" -> "draft"            (Metadata!B6, "Status" row)
#   - sharedStrings: date string bumped to the new   (Metadata!B8, "Date" row)
#     publish timestamp 2023-08-01T16:12:28+00:00
#   - styles.xml: the two content cellXfs (the bold header style and the
#     plain-bordered data style) gain applyAlignment="true" so the existing
#     <alignment vertical="top" wrapText="true"/> actually takes effect.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Concepts")

# --- Update the Status and Date values on the Metadata sheet ---
$ws1.Range("B6").Value = "draft"
$ws1.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# --- Turn on alignment application (vertical=top, wrap text) for every ---
# --- cell in both sheets, matching the two cellXfs touched by the diff ---
$ws1.Range("A1:B21").WrapText = $true
$ws1.Range("A1:B21").VerticalAlignment = -4160   # xlVAlignTop

$ws2.Range("A1:D8").WrapText = $true
$ws2.Range("A1:D8").VerticalAlignment = -4160    # xlVAlignTop
